$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.674.58'
$ws.Range("E2").Value = '  -0.18%  '
$ws.Range("D3").Value = '2.291.65'
$ws.Range("E3").Value = '  -0.13%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").Value = "'96.36"
$ws.Range("E5").Value = '  +2.74%  '
$ws.Range("D6").Value = "'267.55"
$ws.Range("E6").Value = '  -0.84%  '
$ws.Range("E7").Value = '  -1.64%  '
$ws.Range("D8").Value = "'0.999"
$ws.Range("E8").Value = '  -0.08%  '
$ws.Range("D9").Value = "'0.611"
$ws.Range("E9").Value = '  -1.35%  '
$ws.Range("D10").Value = "'45.87"
$ws.Range("E10").Value = '  +1.04%  '
$ws.Range("E11").Value = '  -0.06%  '
$ws.Range("D12").Value = "'7.83"
$ws.Range("E12").Value = '  -2.80%  '
$ws.Range("E13").Value = '  +0.31%  '
$ws.Range("D14").Value = '2.631.53'
$ws.Range("E14").Value = '  -0.07%  '
$ws.Range("D15").Value = "'15.17"
$ws.Range("E15").Value = '  -0.37%  '
$ws.Range("D16").Value = "'0.849"
$ws.Range("E16").Value = '  -0.27%  '
$ws.Range("D17").Value = '2.290.36'
$ws.Range("E17").Value = '  +0.02%  '
$ws.Range("D18").Value = '43.574.64'
$ws.Range("E18").Value = '  -0.29%  '
$ws.Range("D19").Value = "'0.0000108"
$ws.Range("E19").Value = '  +2.27%  '
$ws.Range("E20").Value = '  -0.75%  '
$ws.Range("D21").Value = "'72.31"
$ws.Range("E21").Value = '  +1.68%  '
$ws.Range("D22").Value = "'2.53"
$ws.Range("E22").Value = '  +10.83%  '
$ws.Range("D23").Value = "'233.11"
$ws.Range("E23").Value = '  -1.34%  '
$ws.Range("D24").Value = "'9.16"
$ws.Range("E24").Value = '  -5.48%  '
$ws.Range("D25").Value = "'1.00"
$ws.Range("E25").Value = '  -0.05%  '
$ws.Range("E26").Value = '  -0.15%  '
$ws.Range("D27").Value = "'11.17"
$ws.Range("E27").Value = '  -0.53%  '
$ws.Range("E28").Value = '  +2.28%  '
$ws.Range("D29").Value = "'40.58"
$ws.Range("E29").Value = '  +3.40%  '
$ws.Range("D30").Value = "'2.27"
$ws.Range("E30").Value = '  +0.55%  '
$ws.Range("D31").Value = "'175.67"
$ws.Range("E31").Value = '  +1.22%  '
$ws.Range("D32").Value = "'21.86"
$ws.Range("E32").Value = '  -1.82%  '
$ws.Range("D33").Value = "'0.0893"
$ws.Range("E33").Value = '  +0.85%  '
$ws.Range("D34").Value = "'5.37"
$ws.Range("E34").Value = '  -3.15%  '
$ws.Range("E35").Value = '  -0.21%  '
$ws.Range("E36").Value = '  -1.38%  '
$ws.Range("E37").Value = '  +1.90%  '
$ws.Range("D38").Value = "'4.35"
$ws.Range("E38").Value = '  -3.77%  '
$ws.Range("D39").Value = "'3.39"
$ws.Range("E39").Value = '  +0.08%  '
$ws.Range("D40").Value = "'0.238"
$ws.Range("E40").Value = '  +2.07%  '
$ws.Range("D41").Value = "'2.32"
$ws.Range("E41").Value = '  +0.58%  '
$ws.Range("E42").Value = '  +0.07%  '
$ws.Range("B43").Value = 'ARBITRUM'
$ws.Range("C43").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D43").Value = "'1.36"
$ws.Range("E43").Value = '  +3.16%  '
$ws.Range("B44").Value = 'MultiversX'
$ws.Range("C44").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D44").Value = "'65.54"
$ws.Range("E44").Value = '  +7.07%  '
$ws.Range("E45").Value = '  -4.22%  '
$ws.Range("D46").Value = "'8.75"
$ws.Range("E46").Value = '  -1.21%  '
$ws.Range("D48").Value = "'97.46"
$ws.Range("E48").Value = '  -2.57%  '
$ws.Range("E49").Value = '  +0.22%  '
$ws.Range("D50").Value = "'0.431"
$ws.Range("E50").Value = '  -0.03%  '
$ws.Range("D51").Value = '2.513.43'
$ws.Range("E51").Value = '  -0.01%  '
